$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) columns
# with the latest scraped values from the GitHub Actions run.

$ws.Range("D2").Value = '29.478.70'
$ws.Range("E2").Value = '  +0.55%  '
$ws.Range("D3").Value = '1.913.92'
$ws.Range("E3").Value = '  +0.07%  '
$c = $ws.Range("D4")
$c.Value = '''1.008'
$c.ClearFormats()
$ws.Range("E4").Value = '  +0.70%  '
$c = $ws.Range("D5")
$c.Value = '''325.72'
$c.ClearFormats()
$ws.Range("E5").Value = '  +1.05%  '
$ws.Range("E6").Value = '  +0.57%  '
$c = $ws.Range("D7")
$c.Value = '''0.4821'
$c.ClearFormats()
$ws.Range("E7").Value = '  +1.78%  '
$c = $ws.Range("D8")
$c.Value = '''0.4061'
$c.ClearFormats()
$ws.Range("E8").Value = '  -0.33%  '
$c = $ws.Range("D9")
$c.Value = '''0.08144'
$c.ClearFormats()
$ws.Range("E9").Value = '  +1.27%  '
$c = $ws.Range("D10")
$c.Value = '''1.012'
$c.ClearFormats()
$ws.Range("E10").Value = '  +0.78%  '
$c = $ws.Range("D11")
$c.Value = '''23.52'
$c.ClearFormats()
$ws.Range("E11").Value = '  +4.32%  '
$ws.Range("D12").Value = '1.894.33'
$ws.Range("E12").Value = '  -0.98%  '
$c = $ws.Range("D13")
$c.Value = '''6.002'
$c.ClearFormats()
$ws.Range("E13").Value = '  +1.79%  '
$c = $ws.Range("D14")
$c.Value = '''7.133'
$c.ClearFormats()
$ws.Range("E14").Value = '  -0.07%  '
$ws.Range("E15").Value = '  +0.32%  '
$ws.Range("E16").Value = '  +0.73%  '
$c = $ws.Range("D17")
$c.Value = '''0.06772'
$c.ClearFormats()
$ws.Range("E17").Value = '  +2.09%  '
$c = $ws.Range("D18")
$c.Value = '''0.00001038'
$c.ClearFormats()
$ws.Range("E18").Value = '  +0.82%  '
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("E20").Value = '  +0.38%  '
$ws.Range("D21").Value = '29.488.95'
$ws.Range("E21").Value = '  +0.55%  '
$c = $ws.Range("D22")
$c.Value = '''5.623'
$c.ClearFormats()
$ws.Range("E22").Value = '  +1.82%  '
$ws.Range("E23").Value = '  +2.63%  '
$c = $ws.Range("D24")
$c.Value = '''2.186'
$c.ClearFormats()
$ws.Range("E24").Value = '  -0.57%  '
$ws.Range("D25").Value = '2.132.82'
$ws.Range("E25").Value = '  +0.22%  '
$c = $ws.Range("D26")
$c.Value = '''155.77'
$c.ClearFormats()
$ws.Range("E26").Value = '  +0.52%  '
$c = $ws.Range("D27")
$c.Value = '''6.368'
$c.ClearFormats()
$ws.Range("E27").Value = '  +4.83%  '
$c = $ws.Range("D28")
$c.Value = '''20.02'
$c.ClearFormats()
$ws.Range("E28").Value = '  +1.20%  '
$c = $ws.Range("D29")
$c.Value = '''2.110'
$c.ClearFormats()
$ws.Range("E29").Value = '  -0.05%  '
$c = $ws.Range("D30")
$c.Value = '''119.81'
$c.ClearFormats()
$ws.Range("E30").Value = '  +1.80%  '
$c = $ws.Range("D31")
$c.Value = '''1.023'
$c.ClearFormats()
$ws.Range("E31").Value = '  -4.90%  '
$ws.Range("E32").Value = '  -0.30%  '
$c = $ws.Range("D33")
$c.Value = '''5.524'
$c.ClearFormats()
$ws.Range("E33").Value = '  +2.35%  '
$c = $ws.Range("D34")
$c.Value = '''3.563'
$c.ClearFormats()
$ws.Range("E34").Value = '  +0.41%  '
$ws.Range("E35").Value = '  -2.73%  '
$c = $ws.Range("D36")
$c.Value = '''0.02266'
$c.ClearFormats()
$ws.Range("E36").Value = '  +0.61%  '
$c = $ws.Range("D37")
$c.Value = '''0.06097'
$c.ClearFormats()
$ws.Range("E37").Value = '  +0.13%  '
$c = $ws.Range("D38")
$c.Value = '''1.174'
$c.ClearFormats()
$ws.Range("E38").Value = '  +0.21%  '
$c = $ws.Range("D39")
$c.Value = '''0.5969'
$c.ClearFormats()
$ws.Range("E39").Value = '  +1.58%  '
$c = $ws.Range("D40")
$c.Value = '''7.983'
$c.ClearFormats()
$ws.Range("E40").Value = '  -3.34%  '
$c = $ws.Range("D41")
$c.Value = '''10.68'
$c.ClearFormats()
$ws.Range("E41").Value = '  +5.59%  '
$ws.Range("E42").Value = '  +0.85%  '
$c = $ws.Range("D43")
$c.Value = '''1.283'
$c.ClearFormats()
$ws.Range("E43").Value = '  +0.45%  '
$ws.Range("E44").Value = '  -5.46%  '
$c = $ws.Range("D45")
$c.Value = '''12.54'
$c.ClearFormats()
$ws.Range("E45").Value = '  +3.48%  '
$c = $ws.Range("D46")
$c.Value = '''0.07628'
$c.ClearFormats()
$ws.Range("E46").Value = '  -3.47%  '
$c = $ws.Range("D47")
$c.Value = '''0.5573'
$c.ClearFormats()
$ws.Range("E47").Value = '  +0.50%  '
$ws.Range("E48").Value = '  +0.69%  '
$c = $ws.Range("D49")
$c.Value = '''115.92'
$c.ClearFormats()
$ws.Range("E49").Value = '  +2.54%  '
$ws.Range("E50").Value = '  +3.03%  '
$ws.Range("E51").Value = '  +1.62%  '
